$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 346
$ws.Range("F3").Value = 3518
$ws.Range("F5").Value = 8233
$ws.Range("F7").Value = 91
$ws.Range("F8").Value = 2177
$ws.Range("F10").Value = 187
$ws.Range("F11").Value = 20
$ws.Range("F12").Value = 1195
$ws.Range("F13").Value = 57
$ws.Range("F15").Value = 14
$ws.Range("F16").Value = 582
$ws.Range("F17").Value = 81
$ws.Range("F18").Value = 3191
$ws.Range("F20").Value = 7263
$ws.Range("F22").Value = 55775
$ws.Range("F23").Value = 55775
$ws.Range("F24").Value = 4466
$ws.Range("F26").Value = 1037
$ws.Range("F27").Value = 869
$ws.Range("F28").Value = 431
$ws.Range("F30").Value = 883
$ws.Range("F32").Value = 3323
$ws.Range("F34").Value = 44
$ws.Range("F36").Value = 868
$ws.Range("F37").Value = 1209
$ws.Range("F38").Value = 1172
$ws.Range("F39").Value = 156
$ws.Range("F40").Value = 190
$ws.Range("F41").Value = 1069
$ws.Range("F42").Value = 702
$ws.Range("F43").Value = 4
$ws.Range("F45").Value = 163
$ws.Range("F46").Value = 10
$ws.Range("F47").Value = 163

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 172
$ws.Range("F9").Value = 31
$ws.Range("F16").Value = 7475
$ws.Range("F17").Value = 106
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 119
$ws.Range("F35").Value = 30
$ws.Range("F42").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2290
$ws.Range("F5").Value = 1551
$ws.Range("F7").Value = 659
$ws.Range("F8").Value = 2337
$ws.Range("F9").Value = 9336
$ws.Range("F10").Value = 1668
$ws.Range("F12").Value = 85

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 346
$ws.Range("F3").Value = 3518
$ws.Range("F5").Value = 8233
$ws.Range("F6").Value = 1551
$ws.Range("F7").Value = 659
$ws.Range("F8").Value = 2337
$ws.Range("F9").Value = 1668
$ws.Range("F11").Value = 85
$ws.Range("F13").Value = 91
$ws.Range("F15").Value = 187
$ws.Range("F16").Value = 57
$ws.Range("F17").Value = 14
$ws.Range("F18").Value = 582
$ws.Range("F19").Value = 81
$ws.Range("F20").Value = 7263
$ws.Range("F21").Value = 55775
$ws.Range("F22").Value = 31
$ws.Range("F24").Value = 4466
$ws.Range("F25").Value = 1037
$ws.Range("F26").Value = 431
$ws.Range("F30").Value = 3323
$ws.Range("F32").Value = 44
$ws.Range("F34").Value = 868
$ws.Range("F35").Value = 1209
$ws.Range("F36").Value = 106
$ws.Range("F37").Value = 156
$ws.Range("F38").Value = 190
$ws.Range("F39").Value = 1069
$ws.Range("F40").Value = 703
$ws.Range("F42").Value = 163
$ws.Range("F44").Value = 163
$ws.Range("F47").Value = 30
$ws.Range("F50").Value = 23
